$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.175.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.789.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.552"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.18%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.69"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.17"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.42%  "
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0930"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.046.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +12.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.783.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.633"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.155.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "254.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0744"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.00%  "
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0519"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.456.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.65%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.636"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.43%  "
$ws.Range("E39").Value = "  +1.17%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.903"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0508"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.94%  "
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.944.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.17%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.33%  "
